# Results from R script
#
# 1. Row 109: correct the date/time serial value in A109 (was using the
#    afternoon timestamp of the *next* trading session's scrape; corrected
#    to the 07:00 UTC snapshot time used throughout the sheet).
# 2. Row 110: append the newly scraped observation for BWZ.MI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 109 -------------------------------------------------------------
$ws.Cells.Item(109, 1).Value = 45478.2916666667

# --- Row 110 (new row) -----------------------------------------------------
# Numeric columns
$ws.Cells.Item(110, 2).Value = 14700
$ws.Cells.Item(110, 3).Value = 0.714999973773956
$ws.Cells.Item(110, 4).Value = 0.675000011920929
$ws.Cells.Item(110, 5).Value = 0.675000011920929
$ws.Cells.Item(110, 6).Value = 0.680000007152557

# Text columns (adj_close is stored as text in this sheet, ticker likewise)
# Force text entry so "0.680000007152557" isn't coerced back into a number,
# then drop the temporary text format so the cell keeps the sheet's default
# (unstyled) appearance, matching the rest of column G.
$ws.Cells.Item(110, 7).NumberFormat = "@"
$ws.Cells.Item(110, 7).Value = "0.680000007152557"
$ws.Cells.Item(110, 7).ClearFormats()

$ws.Cells.Item(110, 8).Value = "BWZ.MI"

# Date/time column - copy A109's date-time number format onto A110 so the
# new row reuses the existing style instead of minting a new one.
$ws.Cells.Item(109, 1).Copy()
$ws.Cells.Item(110, 1).PasteSpecial(-4122)
$ws.Cells.Item(110, 1).Value = 45481.331087963
